# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.7287194209349384, 0.3375848360084654, 0.1529057820181812, 0.4998867070740569, 1, 1.719096746035642)
    3 = @(1.505614041169197, 0.3375848360084654, 0.1529057820181812, 0.4998867070740569, 0, 2.495991366269901)
    4 = @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 0, 5.488907176552729)
    5 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 0, 6.048734245549538)
    6 = @(1.505614041169197, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 0, 4.371470058157054)
    7 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 0, 6.048734245549538)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]  # B: TB
    $ws.Cells.Item($row, 3).Value = $vals[1]  # C: d2S
    $ws.Cells.Item($row, 4).Value = $vals[2]  # D: K
    $ws.Cells.Item($row, 5).Value = $vals[3]  # E: IP
    $ws.Cells.Item($row, 6).Value = $vals[4]  # F: Win
    $ws.Cells.Item($row, 7).Value = $vals[5]  # G: sum
}
